# "Generate Report for Handback" — refresh the handback status report:
#   - zh-cn and de-de now disagree with en-US, so the Overview status text
#     (and the per-language sheets' Status column, which share the same
#     string) flips from "in sync" to "not in sync".
#   - The second file's (7e02af04...) handback datetimes were refreshed
#     to the latest generation run for both locales.
#   - The Status column (Overview E/F, and zh-cn/de-de C) grows wider to
#     fit the longer text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Handed back: not in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update handed-back status text (shared across Overview + both language sheets) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Refresh the Correspond Handback DateTime for the second file (row 3) ---
$wsZhCn.Range("K3").Value = "2016-11-03 20:09:32"
$wsDeDe.Range("K3").Value = "2016-11-03 20:09:49"

# --- Widen the Status columns to fit the longer text ---
$wsOverview.Columns.Item(5).ColumnWidth = 32.65
$wsOverview.Columns.Item(6).ColumnWidth = 32.65
$wsZhCn.Columns.Item(3).ColumnWidth = 32.65
$wsDeDe.Columns.Item(3).ColumnWidth = 32.65
